# Progress on generating topology & xyz's
# - Set the LJ sheet's selection (no longer the active/tabSelected sheet)
# - Add a new "Coords" worksheet at the end of the workbook with atom
#   coordinate data and CONCATENATE formulas building python output lines
# - Make "Coords" the active sheet with E1:E9 selected

$wb = $excel.ActiveWorkbook

# --- LJ sheet: update the lingering selection, it's no longer the active tab ---
$ljSheet = $wb.Worksheets.Item("LJ")
$ljSheet.Activate()
$ljSheet.Range("C22").Select()

# --- Add the new "Coords" sheet after the last existing sheet (LJ) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Coords"

# Row data: Element, x-offset, y-offset, z-offset
$ws.Range("A1").Value = "C"
$ws.Range("B1").Value = 2
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 0.5

$ws.Range("A2").Value = "H"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0.2

$ws.Range("A3").Value = "H"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 1.1

$ws.Range("A4").Value = "H"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 0.1

$ws.Range("A5").Value = "C"
$ws.Range("B5").Value = 3.5
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 0.9

$ws.Range("A6").Value = "H"
$ws.Range("B6").Value = 3.5
$ws.Range("C6").Value = 2.8
$ws.Range("D6").Value = 1.7

$ws.Range("A7").Value = "H"
$ws.Range("B7").Value = 3.5
$ws.Range("C7").Value = 1.3
$ws.Range("D7").Value = 1.7

$ws.Range("A8").Value = "O"
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 0.4

$ws.Range("A9").Value = "H"
$ws.Range("B9").Value = 5.5
$ws.Range("C9").Value = 2.5
$ws.Range("D9").Value = 0.2

# Column E: build the python outputFile.write(...) line for each atom
for ($r = 1; $r -le 9; $r++) {
    $a = "A$r"
    $b = "B$r"
    $c = "C$r"
    $d = "D$r"
    $ws.Range("E$r").Formula = "=CONCATENATE(""outputFile.write(f'"",$a,"" {rand[0] + "",$b,""} {rand[1] + "",$c,""} {rand[1] + "",$d,""} \n')"")"
}

# --- Make Coords the active sheet, with E1:E9 selected ---
$ws.Activate()
$ws.Range("E1:E9").Select()
